# Auto-generated Excel COM-interop script
# Refresh of the "cryptos" worksheet: updates Price (D) and Volume(1h) (E)
# columns with the latest scraped figures, and corrects the ordering of two
# Coin/Link row pairs (B/C) that the upstream feed re-ranked, matching commit
# "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric- and percent-looking
# strings (e.g. "67.226.65", "0.619", "  +0.59%  ") are stored verbatim,
# matching the workbook's original inline-string cells instead of being
# auto-converted to numbers by Excel.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

# Coin / Link corrections (row order fix)
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"

# Price / Volume(1h) updates
$ws.Range("D2").Value = "67.226.65"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "3.855.97"
$ws.Range("E3").Value = "  +4.17%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "412.83"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "132.13"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").Value = "3.836.26"
$ws.Range("E7").Value = "  +3.84%  "
$ws.Range("D8").Value = "0.619"
$ws.Range("E8").Value = "  -4.13%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "0.740"
$ws.Range("E10").Value = "  -3.82%  "
$ws.Range("D11").Value = "0.172"
$ws.Range("E11").Value = "  -5.89%  "
$ws.Range("D12").Value = "0.0000379"
$ws.Range("E12").Value = "  -5.12%  "
$ws.Range("D13").Value = "41.08"
$ws.Range("E13").Value = "  -4.99%  "
$ws.Range("D14").Value = "4.452.45"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("E15").Value = "  -6.48%  "
$ws.Range("D16").Value = "15.26"
$ws.Range("E16").Value = "  +15.51%  "
$ws.Range("D17").Value = "3.869.30"
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "19.61"
$ws.Range("E19").Value = "  -5.04%  "
$ws.Range("D20").Value = "67.591.67"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "1.07"
$ws.Range("E21").Value = "  -5.12%  "
$ws.Range("D22").Value = "414.76"
$ws.Range("E22").Value = "  -6.94%  "
$ws.Range("D23").Value = "14.54"
$ws.Range("E23").Value = "  -11.78%  "
$ws.Range("D24").Value = "85.97"
$ws.Range("E24").Value = "  -4.87%  "
$ws.Range("D25").Value = "3.07"
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("D26").Value = "36.94"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  +13.70%  "
$ws.Range("D28").Value = "3.15"
$ws.Range("E28").Value = "  -5.65%  "
$ws.Range("D29").Value = "9.52"
$ws.Range("E29").Value = "  -7.49%  "
$ws.Range("D30").Value = "689.47"
$ws.Range("E30").Value = "  +5.52%  "
$ws.Range("D31").Value = "12.53"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("D33").Value = "2.74"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "7.23"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").Value = "  -8.41%  "
$ws.Range("D36").Value = "39.07"
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("D37").Value = "0.0₃0815"
$ws.Range("E37").Value = "  +9.14%  "
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "54.98"
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("D40").Value = "0.0461"
$ws.Range("E40").Value = "  -7.12%  "
$ws.Range("D41").Value = "3.07"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("D42").Value = "0.995"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "0.138"
$ws.Range("E43").Value = "  -9.29%  "
$ws.Range("D44").Value = "148.87"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "4.49"
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("D46").Value = "3.19"
$ws.Range("E46").Value = "  +18.37%  "
$ws.Range("D47").Value = "3.34"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("D48").Value = "26.97"
$ws.Range("E48").Value = "  -8.53%  "
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").Value = "2.83"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("D51").Value = "2.58"
$ws.Range("E51").Value = "  -3.41%  "

# Restore default (General) formatting now that the text values are set,
# so cell styling matches the rest of the sheet.
$priceVolumeRange.ClearFormats()

